# feat: add 2022-Q1 data
#
# The previously-last sheet "总计" (summary) becomes "2022-Q1" (holding the
# new quarter's fund-holding detail rows), and a fresh "总计" sheet is
# appended after it, carrying the old summary rows plus a new first row for
# 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the existing "总计" sheet right after itself (so the copy
#    inherits its sheetPr/pageMargins instead of generic blank-sheet
#    defaults), then repurpose the original as "2022-Q1" and keep the
#    duplicate as the (still-last) "总计" — giving sheet order/ids
#    [... , 2022-Q1, 总计].
# ---------------------------------------------------------------------
$qSheet = $wb.Worksheets.Item("总计")
$qIndex = $qSheet.Index
$qSheet.Copy($null, $qSheet)
$totalSheet = $wb.Worksheets.Item($qIndex + 1)

$qSheet.Name = "2022-Q1"
$totalSheet.Name = "总计"

# A cell that already carries the workbook's bold/centered/bordered
# "header & index" look, used as a formatting template so the copy always
# lands on the very same style record instead of minting new ones.
$styleTemplate = $wb.Worksheets.Item("2021-Q4").Cells.Item(1, 2)

function Copy-HeaderStyle($cell) {
    $styleTemplate.Copy()
    $cell.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# Writes a text value into a cell while keeping it typed as text (numeric
# looking codes like "009562" or figures like "6.65" would otherwise be
# silently re-interpreted as numbers), then strips the leading-apostrophe
# formatting marker back off so the cell keeps the sheet's plain/default
# style instead of picking up a "quote prefix" style.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 2. Rebuild "2022-Q1" with the per-fund holding detail.
# ---------------------------------------------------------------------
$qSheet.Cells.ClearContents()

$qHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $qHeaders.Length; $i++) {
    $cell = $qSheet.Cells.Item(1, $i + 2)
    Copy-HeaderStyle $cell
    $cell.Value = $qHeaders[$i]
}

$qRows = @(
    @("009562", "工银瑞信中国机会全球配置股票(QDII)美元", "6.65", "92.85", "3.23", "0.2148", 3),
    @("486001", "工银瑞信中国机会全球配置股票(QDII)", "6.65", "92.85", "3.23", "0.2148", 3),
    @("009563", "工银瑞信中国机会全球配置股票(QDII)港币", "6.65", "92.85", "3.23", "0.2148", 3),
    @("486002", "工银全球精选股票(QDII)", "4.23", "94.60", "4.56", "0.1929", 2),
    @("012751", "建信纳斯达克100指数（QDII）A 美元现汇", "0.34", "88.02", "13.39", "0.0455", 1),
    @("012752", "建信纳斯达克100指数（QDII）C 人民币", "0.34", "88.02", "13.39", "0.0455", 1),
    @("012753", "建信纳斯达克100指数（QDII）C 美元现汇", "0.34", "88.02", "13.39", "0.0455", 1)
)

for ($r = 0; $r -lt $qRows.Length; $r++) {
    $row = $r + 2
    $data = $qRows[$r]

    $idxCell = $qSheet.Cells.Item($row, 1)
    Copy-HeaderStyle $idxCell
    $idxCell.Value = $r

    Set-TextValue $qSheet.Cells.Item($row, 2) $data[0]
    Set-TextValue $qSheet.Cells.Item($row, 3) $data[1]
    Set-TextValue $qSheet.Cells.Item($row, 4) $data[2]
    Set-TextValue $qSheet.Cells.Item($row, 5) $data[3]
    Set-TextValue $qSheet.Cells.Item($row, 6) $data[4]
    Set-TextValue $qSheet.Cells.Item($row, 7) $data[5]
    $qSheet.Cells.Item($row, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 3. Rebuild "总计" with the quarterly summary, including the new
#    2022-Q1 row at the top of the data (old rows shift down by one).
# ---------------------------------------------------------------------
$totalSheet.Cells.ClearContents()

$tHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $tHeaders.Length; $i++) {
    $cell = $totalSheet.Cells.Item(1, $i + 2)
    Copy-HeaderStyle $cell
    $cell.Value = $tHeaders[$i]
}

$tRows = @(
    @("2022-Q1", 7, 0.97),
    @("2021-Q4", 3, 2.71),
    @("2021-Q3", 4, 0.06),
    @("2021-Q2", 1, 0.01),
    @("2021-Q1", 1, 0.01),
    @("2020-Q4", 3, 0.51)
)

for ($r = 0; $r -lt $tRows.Length; $r++) {
    $row = $r + 2
    $data = $tRows[$r]

    $idxCell = $totalSheet.Cells.Item($row, 1)
    Copy-HeaderStyle $idxCell
    $idxCell.Value = $r

    Set-TextValue $totalSheet.Cells.Item($row, 2) $data[0]
    $totalSheet.Cells.Item($row, 3).Value = $data[1]
    $totalSheet.Cells.Item($row, 4).Value = $data[2]
}
